# Updates the "scraped_at" timestamps (column K) on the "snapshot" sheet
# to reflect a fresh scrape run (refreshed values, 2025-11-23 ~11:04-11:05 UTC).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$ws.Range("K2").Value = "2025-11-23T11:04:25.502898+00:00"
$ws.Range("K3").Value = "2025-11-23T11:04:25.502935+00:00"
$ws.Range("K4").Value = "2025-11-23T11:04:28.092503+00:00"
$ws.Range("K5").Value = "2025-11-23T11:04:28.092536+00:00"
$ws.Range("K6").Value = "2025-11-23T11:04:30.212849+00:00"
$ws.Range("K7").Value = "2025-11-23T11:04:32.257060+00:00"
$ws.Range("K8").Value = "2025-11-23T11:04:34.346951+00:00"
$ws.Range("K9").Value = "2025-11-23T11:04:34.346968+00:00"
$ws.Range("K10").Value = "2025-11-23T11:04:34.346976+00:00"
$ws.Range("K11").Value = "2025-11-23T11:04:36.430325+00:00"
$ws.Range("K12").Value = "2025-11-23T11:04:38.529464+00:00"
$ws.Range("K13").Value = "2025-11-23T11:04:40.593521+00:00"
$ws.Range("K14").Value = "2025-11-23T11:04:42.642253+00:00"
$ws.Range("K15").Value = "2025-11-23T11:04:45.185838+00:00"
$ws.Range("K16").Value = "2025-11-23T11:04:49.384022+00:00"
$ws.Range("K17").Value = "2025-11-23T11:04:49.384051+00:00"
$ws.Range("K18").Value = "2025-11-23T11:04:51.339839+00:00"
$ws.Range("K19").Value = "2025-11-23T11:04:51.339885+00:00"
$ws.Range("K20").Value = "2025-11-23T11:04:51.339897+00:00"
$ws.Range("K21").Value = "2025-11-23T11:04:53.859383+00:00"
$ws.Range("K22").Value = "2025-11-23T11:04:53.859412+00:00"
$ws.Range("K23").Value = "2025-11-23T11:04:55.902959+00:00"
$ws.Range("K24").Value = "2025-11-23T11:04:55.902978+00:00"
$ws.Range("K25").Value = "2025-11-23T11:04:55.902989+00:00"
$ws.Range("K26").Value = "2025-11-23T11:04:55.902997+00:00"
$ws.Range("K27").Value = "2025-11-23T11:04:58.359933+00:00"
$ws.Range("K28").Value = "2025-11-23T11:04:58.359978+00:00"
$ws.Range("K29").Value = "2025-11-23T11:05:00.396057+00:00"
$ws.Range("K30").Value = "2025-11-23T11:05:00.396087+00:00"
$ws.Range("K31").Value = "2025-11-23T11:05:00.396105+00:00"
$ws.Range("K32").Value = "2025-11-23T11:05:00.396121+00:00"
$ws.Range("K33").Value = "2025-11-23T11:05:02.943085+00:00"
$ws.Range("K34").Value = "2025-11-23T11:05:02.943102+00:00"
$ws.Range("K35").Value = "2025-11-23T11:05:07.979970+00:00"
$ws.Range("K36").Value = "2025-11-23T11:05:07.979999+00:00"
$ws.Range("K37").Value = "2025-11-23T11:05:09.990261+00:00"
$ws.Range("K38").Value = "2025-11-23T11:05:09.990293+00:00"
